# Auto-generated script to apply scheduled-runner price/profit updates
# across all 8 Leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (67 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 103.333336
$ws.Range("J33").Value = 59.666668
$ws.Range("L33").Value = 59.666668
$ws.Range("N33").Value = -517.666668
$ws.Range("H74").Value = 8860
$ws.Range("I74").Value = 4825
$ws.Range("K74").Value = 4825
$ws.Range("M74").Value = -3889
$ws.Range("H77").Value = 8860
$ws.Range("I77").Value = 4825
$ws.Range("K77").Value = 24125
$ws.Range("M77").Value = -19445
$ws.Range("H87").Value = 59536.4
$ws.Range("J87").Value = 59536.4
$ws.Range("L87").Value = 59536.4
$ws.Range("N87").Value = -62032.4
$ws.Range("H90").Value = 59536.4
$ws.Range("J90").Value = 59536.4
$ws.Range("L90").Value = 178609.2
$ws.Range("N90").Value = -191089.2
$ws.Range("H92").Value = 889.1667
$ws.Range("I92").Value = 167.14285
$ws.Range("K92").Value = 167.14285
$ws.Range("M92").Value = 1080.85715
$ws.Range("H108").Value = 34999
$ws.Range("J108").Value = 34999
$ws.Range("L108").Value = 34999
$ws.Range("N108").Value = -42679
$ws.Range("H116").Value = 4999.6665
$ws.Range("I116").Value = 4999.5
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 4999.5
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -1557.5
$ws.Range("N116").Value = -11884
$ws.Range("H117").Value = 93910.664
$ws.Range("J117").Value = 93910.664
$ws.Range("L117").Value = 93910.664
$ws.Range("N117").Value = -103088.664
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 14920.056
$ws.Range("I132").Value = 16404.691
$ws.Range("J132").Value = 11060
$ws.Range("K132").Value = 49214.073
$ws.Range("L132").Value = 33180
$ws.Range("M132").Value = -46684.073
$ws.Range("N132").Value = -38240
$ws.Range("H135").Value = 2261.5
$ws.Range("I135").Value = 2182.8572
$ws.Range("J135").Value = 2445
$ws.Range("K135").Value = 19645.7148
$ws.Range("L135").Value = 22005
$ws.Range("M135").Value = -17110.7148
$ws.Range("N135").Value = -27075
$ws.Range("H138").Value = 2090.375
$ws.Range("I138").Value = 2026.8889
$ws.Range("J138").Value = 2172
$ws.Range("K138").Value = 6080.6667
$ws.Range("L138").Value = 6516
$ws.Range("M138").Value = -940.6666999999998
$ws.Range("N138").Value = -16796
$ws.Range("H140").Value = 72500
$ws.Range("J140").Value = 125000
$ws.Range("L140").Value = 125000
$ws.Range("N140").Value = -135360

# ---- Sheet: ARM (24 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4050.182
$ws.Range("I32").Value = 4050.182
$ws.Range("K32").Value = 4050.182
$ws.Range("M32").Value = -3763.182
$ws.Range("H61").Value = 6429
$ws.Range("I61").Value = 6279.4
$ws.Range("K61").Value = 6279.4
$ws.Range("M61").Value = -6067.4
$ws.Range("H74").Value = 4252.619
$ws.Range("J74").Value = 6171.4287
$ws.Range("L74").Value = 6171.4287
$ws.Range("N74").Value = -7919.4287
$ws.Range("H77").Value = 4252.619
$ws.Range("J77").Value = 6171.4287
$ws.Range("L77").Value = 30857.1435
$ws.Range("N77").Value = -39593.14350000001
$ws.Range("H101").Value = 55554
$ws.Range("J101").Value = 55554
$ws.Range("L101").Value = 55554
$ws.Range("N101").Value = -62044
$ws.Range("H136").Value = 6429
$ws.Range("I136").Value = 6279.4
$ws.Range("K136").Value = 18838.2
$ws.Range("M136").Value = -16288.2

# ---- Sheet: BSM (16 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6310.4614
$ws.Range("J86").Value = 8888
$ws.Range("L86").Value = 8888
$ws.Range("N86").Value = -11134
$ws.Range("H89").Value = 6310.4614
$ws.Range("J89").Value = 8888
$ws.Range("L89").Value = 44440
$ws.Range("N89").Value = -55672
$ws.Range("H105").Value = 2005
$ws.Range("I105").Value = 2005
$ws.Range("K105").Value = 2005
$ws.Range("M105").Value = -258
$ws.Range("H134").Value = 3793.3044
$ws.Range("I134").Value = 3263.7222
$ws.Range("K134").Value = 9791.1666
$ws.Range("M134").Value = -7256.1666

# ---- Sheet: CRP (34 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1765
$ws.Range("I22").Value = 147.5
$ws.Range("K22").Value = 147.5
$ws.Range("M22").Value = 202.5
$ws.Range("H43").Value = 13728.571
$ws.Range("J43").Value = 13728.571
$ws.Range("L43").Value = 13728.571
$ws.Range("N43").Value = -14096.571
$ws.Range("H94").Value = 3148.0715
$ws.Range("I94").Value = 1148.7142
$ws.Range("J94").Value = 5147.4287
$ws.Range("K94").Value = 1148.7142
$ws.Range("L94").Value = 5147.4287
$ws.Range("M94").Value = -697.7141999999999
$ws.Range("N94").Value = -6049.4287
$ws.Range("H101").Value = 13728.571
$ws.Range("J101").Value = 13728.571
$ws.Range("L101").Value = 13728.571
$ws.Range("N101").Value = -20218.571
$ws.Range("H104").Value = 43392.5
$ws.Range("I104").Value = 44000
$ws.Range("J104").Value = 42785
$ws.Range("K104").Value = 44000
$ws.Range("L104").Value = 42785
$ws.Range("M104").Value = -41379
$ws.Range("N104").Value = -48027
$ws.Range("H107").Value = 366.93332
$ws.Range("J107").Value = 292.16666
$ws.Range("L107").Value = 292.16666
$ws.Range("N107").Value = -4132.16666
$ws.Range("H133").Value = 45011.707
$ws.Range("J133").Value = 45312.438
$ws.Range("L133").Value = 45312.438
$ws.Range("N133").Value = -50372.438

# ---- Sheet: CUL (18 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H98").Value = 175.75
$ws.Range("I98").Value = 175.75
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 527.25
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H134").Value = 7066.5557
$ws.Range("J134").Value = 18666.334
$ws.Range("L134").Value = 55999.00199999999
$ws.Range("N134").Value = -66139.00199999999

# ---- Sheet: GSM (31 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 124.5
$ws.Range("I2").Value = 93.75
$ws.Range("K2").Value = 93.75
$ws.Range("M2").Value = 19.25
$ws.Range("H102").Value = 1262.6562
$ws.Range("I102").Value = 876.4828
$ws.Range("J102").Value = 4995.6665
$ws.Range("K102").Value = 876.4828
$ws.Range("L102").Value = 4995.6665
$ws.Range("M102").Value = 745.5172
$ws.Range("N102").Value = -8239.666499999999
$ws.Range("H122").Value = 169023.17
$ws.Range("I122").Value = 239667.38
$ws.Range("J122").Value = 4186.6665
$ws.Range("K122").Value = 719002.14
$ws.Range("L122").Value = 12559.9995
$ws.Range("M122").Value = -716552.14
$ws.Range("N122").Value = -17459.9995
$ws.Range("H132").Value = 14711.267
$ws.Range("I132").Value = 16122.909
$ws.Range("J132").Value = 10829.25
$ws.Range("K132").Value = 48368.727
$ws.Range("L132").Value = 32487.75
$ws.Range("M132").Value = -45838.727
$ws.Range("N132").Value = -37547.75
$ws.Range("H140").Value = 279697
$ws.Range("I140").Value = 279697
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 279697
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ---- Sheet: LTW (37 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H22").Value = 925
$ws.Range("I22").Value = 925
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 925
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 925
$ws.Range("I27").Value = 925
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 925
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H68").Value = 7759.1816
$ws.Range("I68").Value = 6250
$ws.Range("K68").Value = 6250
$ws.Range("M68").Value = -5501
$ws.Range("H71").Value = 7759.1816
$ws.Range("I71").Value = 6250
$ws.Range("K71").Value = 31250
$ws.Range("M71").Value = -27506
$ws.Range("H100").Value = 6440.4814
$ws.Range("I100").Value = 3198.5
$ws.Range("K100").Value = 3198.5
$ws.Range("M100").Value = -2657.5
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490
$ws.Range("H122").Value = 3672
$ws.Range("I122").Value = 3672
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11016
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# ---- Sheet: WVR (20 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 111873.42
$ws.Range("I4").Value = 132799.69
$ws.Range("K4").Value = 132799.69
$ws.Range("M4").Value = -132686.69
$ws.Range("H80").Value = 25000
$ws.Range("J80").Value = 25000
$ws.Range("L80").Value = 25000
$ws.Range("N80").Value = -26996
$ws.Range("H83").Value = 25000
$ws.Range("J83").Value = 25000
$ws.Range("L83").Value = 75000
$ws.Range("N83").Value = -84984
$ws.Range("H96").Value = 775
$ws.Range("H107").Value = 686.3077
$ws.Range("I107").Value = 611.7
$ws.Range("J107").Value = 935
$ws.Range("K107").Value = 1835.1
$ws.Range("L107").Value = 2805
$ws.Range("M107").Value = 84.89999999999986
$ws.Range("N107").Value = -6645
